# Applies two changes described by the upstream diff:
#  1. The cached "datetimeFigureOut" date field text on the slide master and
#     every slide layout is bumped from 3/23/2021 -> 3/29/2021 (the deck was
#     re-saved on a later day, refreshing the auto date placeholder cache).
#  2. A small position fix for the "Rectangle 29" shape on slide 10 (part of
#     a diagram) - it moves straight down a little (y: 4994275 -> 5029200).

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes, [string]$newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePlaceholder = $false
        if ($shp.Type -eq 14) {
            # msoPlaceholder
            if ($shp.PlaceholderFormat.Type -eq 16) {
                # ppPlaceholderDate
                $isDatePlaceholder = $true
            }
        }
        if ($isDatePlaceholder) {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

# 1. Slide master's Date Placeholder.
Set-DatePlaceholderText $p.SlideMaster.Shapes "3/29/2021"

# ... and every custom layout hanging off that master.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes "3/29/2021"
}

# 2. Nudge "Rectangle 29" on slide 10 straight down (218525 -> 5029200 - no
#    change in size, just the vertical offset of the small label box).
$slide10 = $p.Slides.Item(10)
$rect29 = $slide10.Shapes.Item("Rectangle 29")
$rect29.Top = 5029200 / 12700
